$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each Price/Volume cell is stored as text (e.g. "108.00", "1.000",
# "30.790.31"). Force the cell's number format to Text ("@") right before
# writing the new value so Excel doesn't auto-coerce the string into a
# number and silently drop significant trailing zeros / the dotted
# thousands formatting used by some of these price strings.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.790.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.932.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4888"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2953"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06862"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.21"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "105.94"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.938.48"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07772"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.326"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6994"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.65"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.802.06"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007703"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.611"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.01"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.503"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.838"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.70"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.53"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.159"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1035"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.383"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.575"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.549"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04883"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7577"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.146"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.706"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01996"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "79.28"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.652"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.481"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8877"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4439"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "108.00"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.863"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "980.98"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1242"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.10"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.243"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.14%  "

Write-Output "Updated cryptos list."
